$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.942.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.643.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5055"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.71%  "

$ws.Range("E10").Value = "  +0.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.289"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.654.57"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5443"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7873"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.971.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.418"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.975"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.006"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.57%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.869"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.76%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1145"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.884"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.28%  "

$ws.Range("E29").Value = "  +0.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05015"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.67%  "

$ws.Range("E31").Value = "  +0.47%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.205"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.37%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.536"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.373"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8951"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.619"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.145.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5554"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01562"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8246"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("E44").Value = "  +8.93%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.782.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4543"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05074"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.005"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09533"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.46%  "
